$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "test54un"
$ws.Range("B3").Value = "test54pw"
